$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '25.942.95'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -0.36%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.627.28'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -1.04%  '
$ws.Range("E4").Value = '  -0.05%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '214.17'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.91%  '
$ws.Range("E6").Value = '  -0.79%  '
$ws.Range("E7").Value = '  -0.05%  '
$ws.Range("E8").Value = '  -2.04%  '
$ws.Range("E9").Value = '  -3.29%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '18.43'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -5.84%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0789'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -1.06%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.853.92'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -0.98%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.623.43'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -1.76%  '
$ws.Range("E14").Value = '  -2.12%  '
$ws.Range("E15").Value = '  -3.33%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '25.958.85'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -0.39%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.0₃0738'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -3.24%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '61.32'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -3.37%  '
$ws.Range("E19").Value = '  -0.07%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '192.58'
$ws.Range("D20").Style = "Normal"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '4.25'
$ws.Range("D21").Style = "Normal"
$ws.Range("E22").Value = '  -3.57%  '
$ws.Range("E23").Value = '  -2.02%  '
$ws.Range("E24").Value = '  +0.66%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '143.86'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +0.42%  '
$ws.Range("E26").Value = '  +0.02%  '
$ws.Range("E27").Value = '  -3.85%  '
$ws.Range("E28").Value = '  -2.17%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '15.22'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -1.90%  '
$ws.Range("E30").Value = '  -1.42%  '
$ws.Range("E31").Value = '  -2.20%  '
$ws.Range("E32").Value = '  -4.17%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.12'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -5.40%  '
$ws.Range("E34").Value = '  -2.75%  '
$ws.Range("E35").Value = '  -2.69%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.125.30'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -0.52%  '
$ws.Range("E37").Value = '  -5.74%  '
$ws.Range("E38").Value = '  -1.55%  '
$ws.Range("E39").Value = '  -3.49%  '
$ws.Range("E40").Value = '  -2.31%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '98.12'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -1.07%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.764.82'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -0.89%  '
$ws.Range("E43").Value = '  -4.39%  '
$ws.Range("E44").Value = '  -5.40%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0533'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +2.03%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '54.43'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -3.59%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.48'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -0.84%  '
$ws.Range("B48").Value = 'Mantle'
$ws.Range("C48").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.413'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -0.46%  '
$ws.Range("B49").Value = 'USDD'
$ws.Range("C49").Value = 'https://coinranking.com/coin/z2PZIKQL7+usdd-usdd'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.01'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +0.19%  '
$ws.Range("B50").Value = 'EnergySwap'
$ws.Range("C50").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '7.48'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -4.02%  '
$ws.Range("B51").Value = 'Algorand'
$ws.Range("C51").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0926'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -2.66%  '
